$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 3 for SignUpTestCases (added first so its strings register
# ahead of the updated description below, matching shared string order)
$ws.Range("A3").Value = "SignUpTestCases"
$ws.Range("C3").Value = "Y"

# Update existing row 2 description (B2) for LoginTestCases
$ws.Range("B2").Value = "Login Related Test cases"

$ws.Range("B3").Value = "Self Registration Related Test Cases"

# Update selection to match target state
$ws.Range("A3").Select()
